# Auto-generated: applies scheduled-runner price/profit refresh to Sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 50000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 50000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 50000
$ws.Range("N3").Value = -50228
$ws.Range("H12").Value = 896.06665
$ws.Range("I12").Value = 1085.5834
$ws.Range("J12").Value = 138
$ws.Range("K12").Value = 1085.5834
$ws.Range("L12").Value = 138
$ws.Range("M12").Value = -915.5834
$ws.Range("N12").Value = -478
$ws.Range("H33").Value = 1257.909
$ws.Range("I33").Value = 1297.45
$ws.Range("J33").Value = 862.5
$ws.Range("K33").Value = 1297.45
$ws.Range("L33").Value = 862.5
$ws.Range("M33").Value = -1068.45
$ws.Range("H70").Value = 84654.53999999999
$ws.Range("I70").Value = 251850
$ws.Range("J70").Value = 10345.444
$ws.Range("K70").Value = 755550
$ws.Range("L70").Value = 31036.332
$ws.Range("M70").Value = -755280
$ws.Range("N70").Value = -31576.332
$ws.Range("H73").Value = 84654.53999999999
$ws.Range("I73").Value = 251850
$ws.Range("J73").Value = 10345.444
$ws.Range("K73").Value = 755550
$ws.Range("L73").Value = 31036.332
$ws.Range("M73").Value = -754614
$ws.Range("N73").Value = -32908.33199999999
$ws.Range("H76").Value = 3109.7273
$ws.Range("I76").Value = 3075.4443
$ws.Range("J76").Value = 3264
$ws.Range("K76").Value = 3075.4443
$ws.Range("L76").Value = 3264
$ws.Range("M76").Value = -2760.4443
$ws.Range("H79").Value = 3109.7273
$ws.Range("I79").Value = 3075.4443
$ws.Range("J79").Value = 3264
$ws.Range("K79").Value = 3075.4443
$ws.Range("L79").Value = 3264
$ws.Range("M79").Value = -1983.4443
$ws.Range("H86").Value = 5555.9414
$ws.Range("I86").Value = 5069.75
$ws.Range("J86").Value = 5988.1113
$ws.Range("K86").Value = 5069.75
$ws.Range("L86").Value = 5988.1113
$ws.Range("M86").Value = -3946.75
$ws.Range("N86").Value = -8234.1113
$ws.Range("H89").Value = 5555.9414
$ws.Range("I89").Value = 5069.75
$ws.Range("J89").Value = 5988.1113
$ws.Range("K89").Value = 25348.75
$ws.Range("L89").Value = 29940.5565
$ws.Range("M89").Value = -19732.75
$ws.Range("N89").Value = -41172.5565
$ws.Range("H100").Value = 9689.4
$ws.Range("I100").Value = 3155.375
$ws.Range("J100").Value = 12065.409
$ws.Range("K100").Value = 3155.375
$ws.Range("L100").Value = 12065.409
$ws.Range("M100").Value = -2614.375
$ws.Range("N100").Value = -13147.409
$ws.Range("H102").Value = 50000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 50000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 50000
$ws.Range("N102").Value = -56490
$ws.Range("H103").Value = 2487.1667
$ws.Range("I103").Value = 899.6667
$ws.Range("J103").Value = 3016.3333
$ws.Range("K103").Value = 2699.0001
$ws.Range("L103").Value = 9048.999899999999
$ws.Range("M103").Value = -2113.0001
$ws.Range("H116").Value = 2953.2
$ws.Range("I116").Value = 2741.5833
$ws.Range("J116").Value = 3799.6667
$ws.Range("K116").Value = 2741.5833
$ws.Range("L116").Value = 3799.6667
$ws.Range("M116").Value = 700.4167000000002
$ws.Range("N116").Value = -10683.6667
$ws.Range("H138").Value = 358947.4
$ws.Range("I138").Value = 997.8
$ws.Range("J138").Value = 1253821.4
$ws.Range("K138").Value = 2993.4
$ws.Range("L138").Value = 3761464.2
$ws.Range("M138").Value = 2146.6
$ws.Range("N138").Value = -3771744.2
$ws.Range("H141").Value = 3158.95
$ws.Range("I141").Value = 2599.1177
$ws.Range("J141").Value = 6331.3335
$ws.Range("K141").Value = 7797.353099999999
$ws.Range("L141").Value = 18994.0005
$ws.Range("M141").Value = -2617.353099999999
$ws.Range("M3").ClearContents()
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("H63").Value = 205100
$ws.Range("I63").Value = 6666.6665
$ws.Range("J63").Value = 353925
$ws.Range("K63").Value = 6666.6665
$ws.Range("L63").Value = 353925
$ws.Range("M63").Value = -5980.6665
$ws.Range("N63").Value = -355297
$ws.Range("H66").Value = 205100
$ws.Range("I66").Value = 6666.6665
$ws.Range("J66").Value = 353925
$ws.Range("K66").Value = 33333.3325
$ws.Range("L66").Value = 1769625
$ws.Range("M66").Value = -29901.3325
$ws.Range("N66").Value = -1776489
$ws.Range("H97").Value = 1374.7428
$ws.Range("I97").Value = 792.4666999999999
$ws.Range("J97").Value = 4868.4
$ws.Range("K97").Value = 792.4666999999999
$ws.Range("L97").Value = 4868.4
$ws.Range("M97").Value = -296.4666999999999
$ws.Range("H102").Value = 8002128.5
$ws.Range("I102").Value = 1913.3334
$ws.Range("J102").Value = 28574110
$ws.Range("K102").Value = 1913.3334
$ws.Range("L102").Value = 28574110
$ws.Range("M102").Value = -291.3334
$ws.Range("H132").Value = 2903.587
$ws.Range("I132").Value = 1398.8379
$ws.Range("J132").Value = 9089.777
$ws.Range("K132").Value = 4196.5137
$ws.Range("L132").Value = 27269.331
$ws.Range("M132").Value = -1666.5137
$ws.Range("N52").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1480.2903
$ws.Range("I20").Value = 1337.2
$ws.Range("J20").Value = 1740.4546
$ws.Range("K20").Value = 1337.2
$ws.Range("L20").Value = 1740.4546
$ws.Range("M20").Value = -1090.2
$ws.Range("H33").Value = 6478.1665
$ws.Range("I33").Value = 6773.8
$ws.Range("J33").Value = 5000
$ws.Range("K33").Value = 6773.8
$ws.Range("L33").Value = 5000
$ws.Range("M33").Value = -6437.8
$ws.Range("N33").Value = -5672
$ws.Range("H36").Value = 2231.5
$ws.Range("I36").Value = 1097.25
$ws.Range("J36").Value = 4500
$ws.Range("K36").Value = 1097.25
$ws.Range("L36").Value = 4500
$ws.Range("M36").Value = -563.25
$ws.Range("N36").Value = -5568
$ws.Range("H45").Value = 10000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 10000
$ws.Range("N45").Value = -11616
$ws.Range("H134").Value = 2963.75
$ws.Range("I134").Value = 2412.1
$ws.Range("J134").Value = 3423.4583
$ws.Range("K134").Value = 7236.299999999999
$ws.Range("L134").Value = 10270.3749
$ws.Range("M134").Value = -4701.299999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 4500
$ws.Range("I33").Value = 2700
$ws.Range("J33").Value = 9900
$ws.Range("K33").Value = 2700
$ws.Range("L33").Value = 9900
$ws.Range("M33").Value = -2321
$ws.Range("H51").Value = 208333.33
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 208333.33
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 208333.33
$ws.Range("N51").Value = -209805.33
$ws.Range("H58").Value = 2715.238
$ws.Range("I58").Value = 2271.3076
$ws.Range("J58").Value = 3436.625
$ws.Range("K58").Value = 2271.3076
$ws.Range("L58").Value = 3436.625
$ws.Range("M58").Value = -2068.3076
$ws.Range("H61").Value = 208333.33
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 208333.33
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 208333.33
$ws.Range("N61").Value = -209029.33
$ws.Range("H99").Value = 3436
$ws.Range("I99").Value = 3555.7917
$ws.Range("J99").Value = 1998.5
$ws.Range("K99").Value = 3555.7917
$ws.Range("L99").Value = 1998.5
$ws.Range("M99").Value = -2057.7917
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("H126").Value = 3436
$ws.Range("I126").Value = 3555.7917
$ws.Range("J126").Value = 1998.5
$ws.Range("K126").Value = 10667.3751
$ws.Range("L126").Value = 5995.5
$ws.Range("M126").Value = -8197.375100000001
$ws.Range("H132").Value = 1293.5385
$ws.Range("I132").Value = 1262.7826
$ws.Range("J132").Value = 1529.3334
$ws.Range("K132").Value = 3788.3478
$ws.Range("L132").Value = 4588.0002
$ws.Range("M132").Value = -1258.3478
$ws.Range("N132").Value = -9648.0002
$ws.Range("H134").Value = 1460.7894
$ws.Range("I134").Value = 1436.5454
$ws.Range("J134").Value = 1620.8
$ws.Range("K134").Value = 4309.6362
$ws.Range("L134").Value = 4862.4
$ws.Range("M134").Value = -1774.6362
$ws.Range("H136").Value = 2715.238
$ws.Range("I136").Value = 2271.3076
$ws.Range("J136").Value = 3436.625
$ws.Range("K136").Value = 6813.9228
$ws.Range("L136").Value = 10309.875
$ws.Range("M136").Value = -4263.9228
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8533
$ws.Range("I5").Value = 516.6667
$ws.Range("J5").Value = 10937.9
$ws.Range("K5").Value = 1550.0001
$ws.Range("L5").Value = 32813.7
$ws.Range("M5").Value = -1438.0001
$ws.Range("N5").Value = -33037.7
$ws.Range("H22").Value = 4200.4
$ws.Range("I22").Value = 6000
$ws.Range("J22").Value = 3750.5
$ws.Range("K22").Value = 18000
$ws.Range("L22").Value = 11251.5
$ws.Range("M22").Value = -17831
$ws.Range("N22").Value = -11589.5
$ws.Range("H27").Value = 4200.4
$ws.Range("I27").Value = 6000
$ws.Range("J27").Value = 3750.5
$ws.Range("K27").Value = 18000
$ws.Range("L27").Value = 11251.5
$ws.Range("M27").Value = -17898
$ws.Range("N27").Value = -11455.5
$ws.Range("H34").Value = 2036.2
$ws.Range("I34").Value = 2038.7858
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 6116.357400000001
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -6032.357400000001
$ws.Range("H44").Value = 63383.5
$ws.Range("I44").Value = 1093.8
$ws.Range("J44").Value = 167199.67
$ws.Range("K44").Value = 3281.4
$ws.Range("L44").Value = 501599.01
$ws.Range("M44").Value = -2883.4
$ws.Range("N44").Value = -502395.01
$ws.Range("H56").Value = 373172.16
$ws.Range("I56").Value = 373172.16
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 373172.16
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -372642.16
$ws.Range("H129").Value = 57305.5
$ws.Range("I129").Value = 84122.5
$ws.Range("J129").Value = 3671.5
$ws.Range("K129").Value = 252367.5
$ws.Range("L129").Value = 11014.5
$ws.Range("M129").Value = -247367.5
$ws.Range("N129").Value = -21014.5
$ws.Range("H131").Value = 527974.8
$ws.Range("I131").Value = 2001150
$ws.Range("J131").Value = 1840.8572
$ws.Range("K131").Value = 6003450
$ws.Range("L131").Value = 5522.571599999999
$ws.Range("M131").Value = -5998410
$ws.Range("N131").Value = -15602.5716
$ws.Range("H132").Value = 2057.9355
$ws.Range("I132").Value = 1311.5
$ws.Range("J132").Value = 2317.5652
$ws.Range("K132").Value = 11803.5
$ws.Range("L132").Value = 20858.0868
$ws.Range("M132").Value = -9273.5
$ws.Range("N132").Value = -25918.0868
$ws.Range("H134").Value = 2214.44
$ws.Range("I134").Value = 1276.5652
$ws.Range("J134").Value = 13000
$ws.Range("K134").Value = 3829.6956
$ws.Range("L134").Value = 39000
$ws.Range("M134").Value = 1240.3044
$ws.Range("H135").Value = 8533
$ws.Range("I135").Value = 516.6667
$ws.Range("J135").Value = 10937.9
$ws.Range("K135").Value = 4650.0003
$ws.Range("L135").Value = 98441.09999999999
$ws.Range("M135").Value = -2115.0003
$ws.Range("N135").Value = -103511.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 176
$ws.Range("I2").Value = 251.44444
$ws.Range("J2").Value = 40.2
$ws.Range("K2").Value = 251.44444
$ws.Range("L2").Value = 40.2
$ws.Range("M2").Value = -138.44444
$ws.Range("H54").Value = 10000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 10000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 10000
$ws.Range("N54").Value = -10780
$ws.Range("H70").Value = 15753.5
$ws.Range("I70").Value = 2503
$ws.Range("J70").Value = 29004
$ws.Range("K70").Value = 2503
$ws.Range("L70").Value = 29004
$ws.Range("M70").Value = -2233
$ws.Range("H73").Value = 15753.5
$ws.Range("I73").Value = 2503
$ws.Range("J73").Value = 29004
$ws.Range("K73").Value = 2503
$ws.Range("L73").Value = 29004
$ws.Range("M73").Value = -1567
$ws.Range("H97").Value = 2957.3704
$ws.Range("I97").Value = 2615.3333
$ws.Range("J97").Value = 3384.9167
$ws.Range("K97").Value = 2615.3333
$ws.Range("L97").Value = 3384.9167
$ws.Range("M97").Value = -2119.3333
$ws.Range("N97").Value = -4376.9167
$ws.Range("H102").Value = 29866.578
$ws.Range("I102").Value = 3705.4285
$ws.Range("J102").Value = 103117.8
$ws.Range("K102").Value = 3705.4285
$ws.Range("L102").Value = 103117.8
$ws.Range("M102").Value = -2083.4285
$ws.Range("H117").Value = 88996.5
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 88996.5
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 88996.5
$ws.Range("N117").Value = -95880.5
$ws.Range("H126").Value = 5527.9165
$ws.Range("I126").Value = 5452.871
$ws.Range("J126").Value = 5664.7646
$ws.Range("K126").Value = 16358.613
$ws.Range("L126").Value = 16994.2938
$ws.Range("M126").Value = -13888.613
$ws.Range("N126").Value = -21934.2938
$ws.Range("H132").Value = 6104.582
$ws.Range("I132").Value = 6368.06
$ws.Range("J132").Value = 3469.8
$ws.Range("K132").Value = 19104.18
$ws.Range("L132").Value = 10409.4
$ws.Range("M132").Value = -16574.18
$ws.Range("N132").Value = -15469.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 500075
$ws.Range("I2").Value = 500075
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 500075
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -499963
$ws.Range("H32").Value = 7734.5
$ws.Range("I32").Value = 7734.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 7734.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -7417.5
$ws.Range("H35").Value = 849.6667
$ws.Range("I35").Value = 849.6667
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 849.6667
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -513.6667
$ws.Range("H40").Value = 5313.136
$ws.Range("I40").Value = 4694.45
$ws.Range("J40").Value = 11500
$ws.Range("K40").Value = 4694.45
$ws.Range("L40").Value = 11500
$ws.Range("M40").Value = -4558.45
$ws.Range("H42").Value = 24000
$ws.Range("I42").Value = 24000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 24000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -23437
$ws.Range("H43").Value = 28500
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 28500
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 28500
$ws.Range("N43").Value = -28886
$ws.Range("H49").Value = 24000
$ws.Range("I49").Value = 24000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 24000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -23853
$ws.Range("H61").Value = 61851.824
$ws.Range("I61").Value = 65535
$ws.Range("J61").Value = 44663.668
$ws.Range("K61").Value = 65535
$ws.Range("L61").Value = 44663.668
$ws.Range("M61").Value = -65333
$ws.Range("H113").Value = 61851.824
$ws.Range("I113").Value = 65535
$ws.Range("J113").Value = 44663.668
$ws.Range("K113").Value = 65535
$ws.Range("L113").Value = 44663.668
$ws.Range("M113").Value = -63365
$ws.Range("H115").Value = 49999
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 49999
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 49999
$ws.Range("N115").Value = -52349
$ws.Range("H122").Value = 14685.562
$ws.Range("I122").Value = 13719.533
$ws.Range("J122").Value = 18308.166
$ws.Range("K122").Value = 41158.599
$ws.Range("L122").Value = 54924.49800000001
$ws.Range("M122").Value = -38708.599
$ws.Range("N122").Value = -59824.49800000001
$ws.Range("H132").Value = 2673.0908
$ws.Range("I132").Value = 1937.9131
$ws.Range("J132").Value = 6430.6665
$ws.Range("K132").Value = 5813.7393
$ws.Range("L132").Value = 19291.9995
$ws.Range("M132").Value = -3283.7393
$ws.Range("N2").ClearContents()
$ws.Range("N42").ClearContents()
$ws.Range("N49").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 2875
$ws.Range("I38").Value = 2875
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 2875
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -2402
$ws.Range("H62").Value = 263225.66
$ws.Range("I62").Value = 263225.66
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 263225.66
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -262601.66
$ws.Range("H65").Value = 263225.66
$ws.Range("I65").Value = 263225.66
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 1316128.3
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -1313008.3
$ws.Range("H113").Value = 3624311.2
$ws.Range("I113").Value = 7576583
$ws.Range("J113").Value = 1395.6666
$ws.Range("K113").Value = 22729749
$ws.Range("L113").Value = 4186.9998
$ws.Range("M113").Value = -22727579
$ws.Range("N113").Value = -8526.9998
$ws.Range("H132").Value = 1100.8889
$ws.Range("I132").Value = 941.6875
$ws.Range("J132").Value = 2374.5
$ws.Range("K132").Value = 2825.0625
$ws.Range("L132").Value = 7123.5
$ws.Range("M132").Value = -295.0625
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

